$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.792.10"
$ws.Range("E2").Value = "  +5.26%  "

$ws.Range("D3").Value = "3.107.35"
$ws.Range("E3").Value = "  +2.94%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'586.01"
$ws.Range("E5").Value = "  +3.76%  "

$ws.Range("D6").Value = "'144.11"
$ws.Range("E6").Value = "  +3.12%  "

$ws.Range("D8").Value = "3.099.49"
$ws.Range("E8").Value = "  +2.99%  "

$ws.Range("E9").Value = "  +1.59%  "

$ws.Range("E10").Value = "  +11.33%  "

$ws.Range("D11").Value = "'5.71"
$ws.Range("E11").Value = "  +8.25%  "

$ws.Range("D12").Value = "'0.468"
$ws.Range("E12").Value = "  +1.48%  "

$ws.Range("E13").Value = "  +5.48%  "

$ws.Range("D14").Value = "'35.40"
$ws.Range("E14").Value = "  +4.13%  "

$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").Value = "3.617.55"
$ws.Range("E16").Value = "  +2.75%  "

$ws.Range("D17").Value = "'7.19"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").Value = "3.102.32"
$ws.Range("E18").Value = "  +2.64%  "

$ws.Range("D19").Value = "62.720.66"
$ws.Range("E19").Value = "  +5.15%  "

$ws.Range("D20").Value = "'463.40"
$ws.Range("E20").Value = "  +6.41%  "

$ws.Range("D21").Value = "'14.07"
$ws.Range("E21").Value = "  +2.87%  "

$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("D23").Value = "'7.53"
$ws.Range("E23").Value = "  +5.58%  "

$ws.Range("D24").Value = "'13.39"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").Value = "'82.23"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("E28").Value = "  +5.01%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +5.15%  "

$ws.Range("E31").Value = "  +8.46%  "

$ws.Range("D32").Value = "'26.94"
$ws.Range("E32").Value = "  +3.45%  "

$ws.Range("E33").Value = "  +8.41%  "

$ws.Range("D34").Value = "0.0₃0823"
$ws.Range("E34").Value = "  +4.93%  "

$ws.Range("D35").Value = "'2.36"
$ws.Range("E35").Value = "  +11.56%  "

$ws.Range("E36").Value = "  +3.88%  "

$ws.Range("D37").Value = "'6.03"
$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D38").Value = "'3.14"
$ws.Range("E38").Value = "  +12.41%  "

$ws.Range("E39").Value = "  +3.59%  "

$ws.Range("D40").Value = "'8.81"
$ws.Range("E40").Value = "  +1.66%  "

$ws.Range("D41").Value = "'431.79"
$ws.Range("E41").Value = "  +7.17%  "

$ws.Range("D42").Value = "2.905.52"
$ws.Range("E42").Value = "  +4.32%  "

$ws.Range("E43").Value = "  +4.03%  "

$ws.Range("D44").Value = "'0.278"
$ws.Range("E44").Value = "  +9.38%  "

$ws.Range("E45").Value = "  +2.78%  "

$ws.Range("D46").Value = "'2.17"
$ws.Range("E46").Value = "  +7.27%  "

$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'35.07"
$ws.Range("E47").Value = "  +4.20%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'123.69"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("D51").Value = "'24.72"
$ws.Range("E51").Value = "  +4.99%  "
